$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "bổ sung lỗi" — revise the reported bug text on row 10 and drop the
# now-redundant bug note that used to live on row 11.
$ws.Range("D10").Value = "_ nhập số thẻ đúng (copy từ csdl tình trạng =1, =0 luôn) nhưng cứ báo tài khoản không đúng."
$ws.Range("D11").ClearContents()

# Row 10 now wraps onto two lines, so it needs to be taller.
$ws.Rows.Item(10).RowHeight = 33

# Leave the cursor where the author's edit ended up.
$ws.Range("D12").Select()
